$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.76
$ws.Range("H2").Value = 2.74
$ws.Range("I2").Value = 3.1
$ws.Range("J2").Value = 2.92
$ws.Range("L2").Value = 1.56
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 1.56
$ws.Range("Q2").Value = 2.5
$ws.Range("T2").Value = 2
$ws.Range("U2").Value = 1.81
$ws.Range("V2").Value = 1.48
$ws.Range("W2").Value = 1.47
$ws.Range("X2").Value = 9.4
$ws.Range("Y2").Value = 9.199999999999999
$ws.Range("Z2").Value = 18.5
$ws.Range("AB2").Value = 9.199999999999999
$ws.Range("AH2").Value = 22
$ws.Range("AN2").Value = 55
$ws.Range("F3").Value = 2.58
$ws.Range("H3").Value = 2.96
$ws.Range("I3").Value = 3.3
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 3.45
$ws.Range("L3").Value = 1.51
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 2.9
$ws.Range("P3").Value = 1.64
$ws.Range("Q3").Value = 2.32
$ws.Range("R3").Value = 1.23
$ws.Range("S3").Value = 4.5
$ws.Range("T3").Value = 1.91
$ws.Range("U3").Value = 1.92
$ws.Range("V3").Value = 1.43
$ws.Range("X3").Value = 10.5
$ws.Range("Z3").Value = 38
$ws.Range("AI3").Value = 160
$ws.Range("AK3").Value = 110
$ws.Range("AL3").Value = 150
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000
$ws.Range("G4").Value = 1.32
$ws.Range("I4").Value = 13.5
$ws.Range("J4").Value = 6.2
$ws.Range("K4").Value = 6.4
$ws.Range("L4").Value = 1.37
$ws.Range("N4").Value = 4.6
$ws.Range("O4").Value = 1.26
$ws.Range("P4").Value = 2.22
$ws.Range("Q4").Value = 1.78
$ws.Range("R4").Value = 1.46
$ws.Range("S4").Value = 3.05
$ws.Range("U4").Value = 1.71
$ws.Range("V4").Value = 1.08
$ws.Range("W4").Value = 4.1
$ws.Range("X4").Value = 21
$ws.Range("Z4").Value = 1000
$ws.Range("AB4").Value = 7.6
$ws.Range("AC4").Value = 14
$ws.Range("AE4").Value = 250
$ws.Range("AH4").Value = 36
$ws.Range("AI4").Value = 210
$ws.Range("AL4").Value = 44
$ws.Range("AM4").Value = 240
$ws.Range("AN4").Value = 5.5
$ws.Range("AO4").Value = 360
$ws.Range("F5").Value = 1.93
$ws.Range("I5").Value = 5.2
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 3.5
$ws.Range("P5").Value = 1.58
$ws.Range("Q5").Value = 2.52
$ws.Range("F6").Value = 2.58
$ws.Range("G6").Value = 2.88
$ws.Range("H6").Value = 2.74
$ws.Range("I6").Value = 3.1
$ws.Range("J6").Value = 3.25
$ws.Range("K6").Value = 3.6
$ws.Range("L6").Value = 1.43
$ws.Range("N6").Value = 3.75
$ws.Range("O6").Value = 1.31
$ws.Range("P6").Value = 1.92
$ws.Range("Q6").Value = 1.94
$ws.Range("S6").Value = 3.4
$ws.Range("T6").Value = 1.7
$ws.Range("V6").Value = 1.48
$ws.Range("W6").Value = 1.53
$ws.Range("F7").Value = 1.98
$ws.Range("G7").Value = 2.1
$ws.Range("H7").Value = 4.7
$ws.Range("I7").Value = 5.5
$ws.Range("J7").Value = 3.1
$ws.Range("K7").Value = 3.5
$ws.Range("L7").Value = 1.56
$ws.Range("M7").Value = 1.12
$ws.Range("N7").Value = 2.72
$ws.Range("O7").Value = 1.52
$ws.Range("P7").Value = 1.55
$ws.Range("Q7").Value = 2.56
$ws.Range("T7").Value = 2.18
$ws.Range("U7").Value = 1.73
$ws.Range("V7").Value = 1.22
$ws.Range("W7").Value = 1.9
$ws.Range("Y7").Value = 1000
$ws.Range("AD7").Value = 46
$ws.Range("AF7").Value = 34
$ws.Range("AG7").Value = 40
$ws.Range("AJ7").Value = 1000
$ws.Range("AK7").Value = 34
$ws.Range("F8").Value = 4.6
$ws.Range("G8").Value = 4.7
$ws.Range("H8").Value = 2.04
$ws.Range("I8").Value = 2.08
$ws.Range("J8").Value = 3.35
$ws.Range("K8").Value = 3.45
$ws.Range("L8").Value = 1.57
$ws.Range("M8").Value = 1.12
$ws.Range("N8").Value = 2.84
$ws.Range("O8").Value = 1.53
$ws.Range("Q8").Value = 2.6
$ws.Range("T8").Value = 2.22
$ws.Range("U8").Value = 1.75
$ws.Range("V8").Value = 1.93
$ws.Range("W8").Value = 1.27
$ws.Range("X8").Value = 9
$ws.Range("Y8").Value = 6.8
$ws.Range("Z8").Value = 10.5
$ws.Range("AA8").Value = 25
$ws.Range("AE8").Value = 27
$ws.Range("AF8").Value = 30
$ws.Range("AG8").Value = 19.5
$ws.Range("AH8").Value = 25
$ws.Range("AJ8").Value = 120
$ws.Range("AK8").Value = 80
$ws.Range("AL8").Value = 110
$ws.Range("AN8").Value = 120
$ws.Range("AO8").Value = 25
$ws.Range("G9").Value = 2.28
$ws.Range("H9").Value = 3.6
$ws.Range("N9").Value = 3.7
$ws.Range("O9").Value = 1.35
$ws.Range("P9").Value = 1.9
$ws.Range("Q9").Value = 2.06
$ws.Range("R9").Value = 1.34
$ws.Range("S9").Value = 3.75
$ws.Range("T9").Value = 1.82
$ws.Range("U9").Value = 2.12
$ws.Range("V9").Value = 1.37
$ws.Range("W9").Value = 1.78
$ws.Range("Y9").Value = 13.5
$ws.Range("AB9").Value = 9.4
$ws.Range("AE9").Value = 44
$ws.Range("AH9").Value = 18
$ws.Range("AI9").Value = 55
$ws.Range("AK9").Value = 24
$ws.Range("AL9").Value = 40
$ws.Range("AM9").Value = 100
$ws.Range("AN9").Value = 18.5
$ws.Range("AO9").Value = 44
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 3.55
$ws.Range("J10").Value = 3.5
$ws.Range("K10").Value = 3.6
$ws.Range("L10").Value = 1.44
$ws.Range("N10").Value = 3.7
$ws.Range("O10").Value = 1.36
$ws.Range("P10").Value = 1.9
$ws.Range("Q10").Value = 2.06
$ws.Range("R10").Value = 1.35
$ws.Range("S10").Value = 3.75
$ws.Range("T10").Value = 1.81
$ws.Range("U10").Value = 2.12
$ws.Range("V10").Value = 1.39
$ws.Range("Z10").Value = 25
$ws.Range("AB10").Value = 10
$ws.Range("AC10").Value = 7.6
$ws.Range("AD10").Value = 14.5
$ws.Range("AE10").Value = 42
$ws.Range("AG10").Value = 11
$ws.Range("AH10").Value = 18
$ws.Range("AI10").Value = 55
$ws.Range("AJ10").Value = 30
$ws.Range("AK10").Value = 25
$ws.Range("AL10").Value = 40
$ws.Range("AM10").Value = 100
$ws.Range("AN10").Value = 20
$ws.Range("AO10").Value = 48
$ws.Range("F11").Value = 2.48
$ws.Range("G11").Value = 2.72
$ws.Range("H11").Value = 2.9
$ws.Range("I11").Value = 3.2
$ws.Range("J11").Value = 3.35
$ws.Range("K11").Value = 3.85
$ws.Range("N11").Value = 3.9
$ws.Range("P11").Value = 1.99
$ws.Range("Q11").Value = 1.92
$ws.Range("S11").Value = 3.25
$ws.Range("V11").Value = 1.46
$ws.Range("W11").Value = 1.58
$ws.Range("AC11").Value = 8.6
$ws.Range("AM11").Value = 100
